# Update cryptocurrency price/volume table with latest values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.388.08"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "1.709.53"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.50"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5328"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06626"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.97"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07624"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.554"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "1.711.30"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "1.943.81"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5777"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "0.0₅8190"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.82"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "27.344.25"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.80"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.661"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.966"
$ws.Range("E23").Value = "  -3.75%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.52"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.727"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.262"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.27"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05407"
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.292"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.505"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.428"
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.648"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.875"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9491"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.411"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5869"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01639"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.851"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "1.047.43"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.90"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "1.850.87"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.03"
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4518"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.087"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05233"
$ws.Range("E51").Value = "  -1.73%  "
